$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (Venezuela exchange-rate entry) appended under the header row
$ws.Range("A2").Value = "VENEZUELA"
$ws.Range("B2").Value = "VES"
$ws.Range("C2").Value = 137
$ws.Range("D2").Value = "Bs."
$ws.Range("E2").Value = '[{"name": "Pago Móvil", "details": "C.I.: V-32147818, Teléfono: 04126027407, Banco: 0105 (Banco Mercantil)"}, {"name": "¡IMPORTANTE!", "details": "¡No colocar conceptos EN LOS PAGOS!"}]'

# Payment-methods cell gets its own readable font/size
$ws.Range("E2").Font.Name = "Arial Unicode MS"
$ws.Range("E2").Font.Size = 10

# A and B columns get an explicit (near-default) width, like C/D/E already have
$ws.Columns("A:B").ColumnWidth = 8.25

# Sheet was left scrolled/selected further down, on an otherwise-empty styled cell
$ws.Range("E12").Font.Underline = $true
$ws.Range("E12").Select()

# Print orientation recorded for the sheet
$ws.PageSetup.Orientation = 1
